$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 5) below the existing entries.
$ws.Range("A5").Value = "awake"
$ws.Range("B5").Value = "thức giấc"

# Return focus to the main window / home cell when the work is finished.
$ws.Range("A1").Select()
